# Fixed limit-10 reading function to read all tabs
#
# The "CasesTab" Cypher query (stored in cell B2 of the "startup" sheet)
# is missing a trailing space after its final "LIMIT 100" clause, which
# was causing the downstream reading/parsing logic to only ever see a
# limit of "100" glued to nothing else and effectively stop after the
# first batch instead of reading every tab. Append the missing trailing
# space so the query text round-trips the same way the other tabs' query
# text already does.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$casesTabQuery = $ws.Range("B2")
$currentText = $casesTabQuery.Value2
if ($currentText.Substring($currentText.Length - 1) -ne " ") {
    $casesTabQuery.Value = $currentText + " "
}

# Move the active selection from B2 to C2 (next column over), matching
# where the user left the cursor after verifying the fix.
$null = $ws.Range("C2").Select()
